$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-4 with new values
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 225

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 175

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 56

# Remove row 5 entirely (it no longer exists in the data)
$ws.Range("A5:B5").ClearContents()
$ws.Rows("5:5").Delete()
